# Adds two new columns (IP Address - 3 / IP Address - 4) with IP data,
# extends the title bar merge/format over the new columns, and normalizes
# the plain data-cell borders across the whole table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New header cells for columns G and H
# ---------------------------------------------------------------------
$ws.Range("G2").Value2 = "IP Address - 3"
$ws.Range("H2").Value2 = "IP Address - 4"

# ---------------------------------------------------------------------
# 2. New IP data for columns G and H, rows 3-17
# ---------------------------------------------------------------------
$gValues = @(
    "13.233.89.201", "3.109.184.3", "13.233.129.132", "65.0.17.121", "13.201.39.151",
    "13.233.63.29", "65.2.172.63", "13.233.84.215", "13.127.77.159", "15.206.171.207",
    "13.233.208.121", "13.233.214.159", "13.126.198.127", "65.2.146.151", "13.233.252.249"
)
$hValues = @(
    "65.0.95.87", "13.232.214.56", "13.233.142.152", "3.110.160.62", "65.0.86.120",
    "15.206.122.37", "13.201.52.146", "13.126.86.253", "35.154.145.87", "13.201.35.58",
    "65.1.131.7", "13.127.0.17", "13.235.246.178", "13.233.247.245", "52.66.188.104"
)

for ($i = 0; $i -lt 15; $i++) {
    $row = $i + 3
    $ws.Cells.Item($row, 7).Value2 = $gValues[$i]
}
# Rows 12/13 in column H are written out of sequence (13 before 12) to
# mirror the shared-string insertion order of the original edit.
$hOrder = @(3, 4, 5, 6, 7, 8, 9, 10, 11, 13, 12, 14, 15, 16, 17)
foreach ($row in $hOrder) {
    $ws.Cells.Item($row, 8).Value2 = $hValues[$row - 3]
}

# ---------------------------------------------------------------------
# 3. Column widths for the two new columns (match existing IP columns)
# ---------------------------------------------------------------------
$ws.Range("G1:H1").ColumnWidth = 13.77734375

# ---------------------------------------------------------------------
# 4. Extend the title merge across the new columns and re-style it:
#    bold, yellow fill, centered - applied uniformly across A1:H1
#    (previously A1 and B1:F1 used two different styles).
# ---------------------------------------------------------------------
$ws.Range("A1:F1").UnMerge()
$ws.Range("A1:H1").Merge()
$ws.Range("A1:H1").Font.Bold = $true
$ws.Range("A1:H1").Interior.Color = 65535
$ws.Range("A1:H1").HorizontalAlignment = -4108

# ---------------------------------------------------------------------
# 5. Normalize all cell borders (title row, header row 2, and body rows
#    3-17, including the two new columns) to the plain thin box-border
#    style used throughout the table - one single pass over the whole
#    used range keeps the border/style table from fragmenting.
# ---------------------------------------------------------------------
$ws.Range("A1:H17").Borders.LineStyle = 1

# ---------------------------------------------------------------------
# 6. Dimension / selection bookkeeping to mirror the saved workbook state
# ---------------------------------------------------------------------
$ws.Range("I14").Select()
